$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewUsers")

$ws.Range("A2").Value = "firstName40"
$ws.Range("B2").Value = "lastName40"
$ws.Range("C2").Value = "fl40@test.com"

$ws.Activate()
$ws.Range("C11").Select()
